$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9927999675585966
$ws.Range("D2").Value = 0.2798731287734992
$ws.Range("E2").Value = 0.2308497643358756
$ws.Range("F2").Value = 0.9456097004347441
$ws.Range("G2").Value = 0.447322835736145
$ws.Range("H2").Value = 0.5737145539869601
$ws.Range("I2").Value = 0.7240609942466598
$ws.Range("J2").Value = 0.2588174211701926
$ws.Range("L2").Value = 0.3734397509095402
$ws.Range("N2").Value = 1.396665933390068
$ws.Range("O2").Value = 1.995591510310362
$ws.Range("B3").Value = 0.9294340059196031
$ws.Range("D3").Value = 0.2819548258911855
$ws.Range("E3").Value = 0.2318011573563368
$ws.Range("F3").Value = 0.9438628670045546
$ws.Range("G3").Value = 0.4387955422283625
$ws.Range("H3").Value = 0.5733800915866851
$ws.Range("I3").Value = 0.7353335356100101
$ws.Range("J3").Value = 0.2580107726554104
$ws.Range("L3").Value = 0.3391230218221608
$ws.Range("N3").Value = 1.383195247513143
$ws.Range("O3").Value = 1.976524758828361
$ws.Range("B4").Value = 0.8906315630160009
$ws.Range("D4").Value = 0.2833116060989109
$ws.Range("E4").Value = 0.2324406609767831
$ws.Range("F4").Value = 0.9433451313041274
$ws.Range("G4").Value = 0.4338788003849459
$ws.Range("H4").Value = 0.5734558444510327
$ws.Range("I4").Value = 0.7426503674489262
$ws.Range("J4").Value = 0.25756497529256
$ws.Range("L4").Value = 0.3180450271966038
$ws.Range("N4").Value = 1.375397648533379
$ws.Range("O4").Value = 1.966071819552155
$ws.Range("B5").Value = 0.8748470234240244
$ws.Range("D5").Value = 0.2838843084449731
$ws.Range("E5").Value = 0.2327152342278378
$ws.Range("F5").Value = 0.9432737251244703
$ws.Range("G5").Value = 0.4319553854449794
$ws.Range("H5").Value = 0.5735574789222113
$ws.Range("I5").Value = 0.7457314432994799
$ws.Range("J5").Value = 0.2573958593714778
$ws.Range("L5").Value = 0.3094544261978314
$ws.Range("N5").Value = 1.372339894895632
$ws.Range("O5").Value = 1.962127608151036
$ws.Range("B6").Value = 0.8722277338304423
$ws.Range("D6").Value = 0.2839806023770231
$ws.Range("E6").Value = 0.2327616722542221
$ws.Range("F6").Value = 0.9432702999396554
$ws.Range("G6").Value = 0.4316408466794002
$ws.Range("H6").Value = 0.5735786314959483
$ws.Range("I6").Value = 0.7462490566679338
$ws.Range("J6").Value = 0.2573685389140365
$ws.Range("L6").Value = 0.3080279150616008
$ws.Range("N6").Value = 1.371839417710348
$ws.Range("O6").Value = 1.961491729591046
$ws.Range("B7").Value = 0.89041857277644
$ws.Range("D7").Value = 0.2833192495302077
$ws.Range("E7").Value = 0.2324443073308302
$ws.Range("F7").Value = 0.9433436030820701
$ws.Range("G7").Value = 0.4338525358641903
$ws.Range("H7").Value = 0.5734569285056494
$ws.Range("I7").Value = 0.7426915174583375
$ws.Range("J7").Value = 0.2575626435799734
$ws.Range("L7").Value = 0.3179291750505513
$ws.Range("N7").Value = 1.375355924370851
$ws.Range("O7").Value = 1.966017349196562
$ws.Range("B8").Value = 0.9709305649010957
$ws.Range("D8").Value = 0.2805746044530726
$ws.Range("E8").Value = 0.2311663499520424
$ws.Range("F8").Value = 0.9448922705266725
$ws.Range("G8").Value = 0.4443163837538435
$ws.Range("H8").Value = 0.5735409302397727
$ws.Range("I8").Value = 0.7278656746287679
$ws.Range("J8").Value = 0.2585290804543305
$ws.Range("L8").Value = 0.3616092902963999
$ws.Range("N8").Value = 1.391923431951057
$ws.Range("O8").Value = 1.988757041875687
$ws.Range("B9").Value = 1.129586050619935
$ws.Range("D9").Value = 0.2758144296071663
$ws.Range("E9").Value = 0.2290970707473754
$ws.Range("F9").Value = 0.9523304038490465
$ws.Range("G9").Value = 0.4673706183385775
$ws.Range("H9").Value = 0.5759335330992883
$ws.Range("I9").Value = 0.701930862926011
$ws.Range("J9").Value = 0.260812476012191
$ws.Range("L9").Value = 0.4471801951866894
$ws.Range("N9").Value = 1.4281378078811
$ws.Range("O9").Value = 2.043301883918474
$ws.Range("B10").Value = 1.246554971967669
$ws.Range("D10").Value = 0.2726940099349813
$ws.Range("E10").Value = 0.2278399150169239
$ws.Range("F10").Value = 0.9604790488143635
$ws.Range("G10").Value = 0.4858606267062555
$ws.Range("H10").Value = 0.5790470482687056
$ws.Range("I10").Value = 0.6847905990310901
$ws.Range("J10").Value = 0.2627210711812964
$ws.Range("L10").Value = 0.5099665630676498
$ws.Range("N10").Value = 1.456978274113069
$ws.Range("O10").Value = 2.089453006822907
$ws.Range("B11").Value = 1.299841844120579
$ws.Range("D11").Value = 0.2713557857629691
$ws.Range("E11").Value = 0.2273244802582877
$ws.Range("F11").Value = 0.9647691466894912
$ws.Range("G11").Value = 0.494610831130359
$ws.Range("H11").Value = 0.5807573761859572
$ws.Range("I11").Value = 0.6774087853666346
$ws.Range("J11").Value = 0.2636383133793387
$ws.Range("L11").Value = 0.538505707629497
$ws.Range("N11").Value = 1.470576038623847
$ws.Range("O11").Value = 2.111770454884123
$ws.Range("B12").Value = 1.320029924964103
$ws.Range("D12").Value = 0.2708606842800876
$ws.Range("E12").Value = 0.2271373623671646
$ws.Range("F12").Value = 0.9664775417684126
$ws.Range("G12").Value = 0.4979731264065919
$ws.Range("H12").Value = 0.5814472473576444
$ws.Range("I12").Value = 0.6746732473708486
$ws.Range("J12").Value = 0.2639925933031293
$ws.Range("L12").Value = 0.5493088312670125
$ws.Range("N12").Value = 1.475793205727939
$ws.Range("O12").Value = 2.120411784877319
$ws.Range("B13").Value = 1.315681663374789
$ws.Range("D13").Value = 0.2709667953869577
$ws.Range("E13").Value = 0.2271773036323488
$ws.Range("F13").Value = 0.9661058807106429
$ws.Range("G13").Value = 0.4972468258131784
$ws.Range("H13").Value = 0.5812967951391386
$ws.Range("I13").Value = 0.6752597340259876
$ws.Range("J13").Value = 0.2639159858268201
$ws.Range("L13").Value = 0.5469823768006563
$ws.Range("N13").Value = 1.474666583651242
$ws.Range("O13").Value = 2.1185422617923
$ws.Range("B14").Value = 1.301502548411861
$ws.Range("D14").Value = 0.2713148200637967
$ws.Range("E14").Value = 0.2273089245957358
$ws.Range("F14").Value = 0.9649080175802283
$ws.Range("G14").Value = 0.4948864713197736
$ws.Range("H14").Value = 0.5808132867882847
$ws.Range("I14").Value = 0.6771825322705731
$ws.Range("J14").Value = 0.2636673217253787
$ws.Range("L14").Value = 0.539394572297482
$ws.Range("N14").Value = 1.471003899832823
$ws.Range("O14").Value = 2.112477570833647
$ws.Range("B15").Value = 1.292818620950698
$ws.Range("D15").Value = 0.271529512112755
$ws.Range("E15").Value = 0.2273905951479911
$ws.Range("F15").Value = 0.9641852074810089
$ws.Range("G15").Value = 0.4934470399687996
$ws.Range("H15").Value = 0.5805226186091232
$ws.Range("I15").Value = 0.6783680906336871
$ws.Range("J15").Value = 0.263515908411911
$ws.Range("L15").Value = 0.5347462752048955
$ws.Range("N15").Value = 1.468769231948556
$ws.Range("O15").Value = 2.10878753989914
$ws.Range("B16").Value = 1.24307396533834
$ws.Range("D16").Value = 0.2727830986935151
$ws.Range("E16").Value = 0.2278747314254996
$ws.Range("F16").Value = 0.960210413621283
$ws.Range("G16").Value = 0.4852956049048203
$ws.Range("H16").Value = 0.5789411841918906
$ws.Range("I16").Value = 0.68528138183567
$ws.Range("J16").Value = 0.2626621042872159
$ws.Range("L16").Value = 0.5081009405782879
$ws.Range("N16").Value = 1.456099183495624
$ws.Range("O16").Value = 2.088021127723607
$ws.Range("B17").Value = 1.212575849845791
$ws.Range("D17").Value = 0.2735729262078186
$ws.Range("E17").Value = 0.2281861551849182
$ws.Range("F17").Value = 0.9579213591224942
$ws.Range("G17").Value = 0.4803818161322511
$ws.Range("H17").Value = 0.5780462723640909
$ws.Range("I17").Value = 0.6896288957550398
$ws.Range("J17").Value = 0.2621507996431376
$ws.Range("L17").Value = 0.4917485324924087
$ws.Range("N17").Value = 1.448448426638706
$ws.Range("O17").Value = 2.07562043607129
$ws.Range("B18").Value = 1.195041494513646
$ws.Range("D18").Value = 0.2740348655407372
$ws.Range("E18").Value = 0.2283705944679255
$ws.Range("F18").Value = 0.9566596552188571
$ws.Range("G18").Value = 0.4775874418199777
$ws.Range("H18").Value = 0.5775592169595711
$ws.Range("I18").Value = 0.6921685643369635
$ws.Range("J18").Value = 0.2618613304589488
$ws.Range("L18").Value = 0.4823409701092203
$ws.Range("N18").Value = 1.44409295887597
$ws.Range("O18").Value = 2.068612423739353
$ws.Range("B19").Value = 1.189105974706592
$ws.Range("D19").Value = 0.2741925853678557
$ws.Range("E19").Value = 0.2284339571427338
$ws.Range("F19").Value = 0.9562418942524573
$ws.Range("G19").Value = 0.4766467936867116
$ws.Range("H19").Value = 0.5773990634487944
$ws.Range("I19").Value = 0.6930351648978004
$ws.Range("J19").Value = 0.2617641178983874
$ws.Range("L19").Value = 0.4791553956462167
$ws.Range("N19").Value = 1.442626032177003
$ws.Range("O19").Value = 2.066261022345714
$ws.Range("B20").Value = 1.215821676725625
$ws.Range("D20").Value = 0.273488056055597
$ws.Range("E20").Value = 0.2281524537236752
$ws.Range("F20").Value = 0.9581593512590274
$ws.Range("G20").Value = 0.4809015948254398
$ws.Range("H20").Value = 0.5781386735559408
$ws.Range("I20").Value = 0.6891620489780306
$ws.Range("J20").Value = 0.2622047515532842
$ws.Range("L20").Value = 0.4934894961006933
$ws.Range("N20").Value = 1.449258206285577
$ws.Range("O20").Value = 2.076927622112805
$ws.Range("B21").Value = 1.305667053050001
$ws.Range("D21").Value = 0.2712122807182702
$ws.Range("E21").Value = 0.2272700458278507
$ws.Range("F21").Value = 0.9652575842807067
$ws.Range("G21").Value = 0.4955784404390471
$ws.Range("H21").Value = 0.5809541599501955
$ws.Range("I21").Value = 0.67661613649959
$ws.Range("J21").Value = 0.2637401729015494
$ws.Range("L21").Value = 0.5416234099295991
$ws.Range("N21").Value = 1.472077879680597
$ws.Range("O21").Value = 2.11425375618893
$ws.Range("B22").Value = 1.364440831109164
$ws.Range("D22").Value = 0.269792849489412
$ws.Range("E22").Value = 0.2267403326609649
$ws.Range("F22").Value = 0.970385284232222
$ws.Range("G22").Value = 0.5054549651185738
$ws.Range("H22").Value = 0.5830402275581861
$ws.Range("I22").Value = 0.6687652159247115
$ws.Range("J22").Value = 0.2647840626410343
$ws.Range("L22").Value = 0.5730578535600159
$ws.Range("N22").Value = 1.487387798566346
$ws.Range("O22").Value = 2.139757123904985
$ws.Range("B23").Value = 1.333067685671836
$ws.Range("D23").Value = 0.2705442224136494
$ws.Range("E23").Value = 0.2270187678288274
$ws.Range("F23").Value = 0.9676038417672572
$ws.Range("G23").Value = 0.5001576485702373
$ws.Range("H23").Value = 0.5819043683092389
$ws.Range("I23").Value = 0.6729234857957636
$ws.Range("J23").Value = 0.2642232576539101
$ws.Range("L23").Value = 0.5562831353635715
$ws.Range("N23").Value = 1.47918062716451
$ws.Range("O23").Value = 2.126044082148724
$ws.Range("B24").Value = 1.214354240082741
$ws.Range("D24").Value = 0.2735264013763441
$ws.Range("E24").Value = 0.2281676733378788
$ws.Range("F24").Value = 0.9580515858581435
$ws.Range("G24").Value = 0.480666507476144
$ws.Range("H24").Value = 0.5780968134931044
$ws.Range("I24").Value = 0.6893729851136317
$ws.Range("J24").Value = 0.2621803459110978
$ws.Range("L24").Value = 0.4927024261599229
$ws.Range("N24").Value = 1.44889197073941
$ws.Range("O24").Value = 2.076336265424175
$ws.Range("B25").Value = 1.086589730606988
$ws.Range("D25").Value = 0.277035829836771
$ws.Range("E25").Value = 0.2296104290391359
$ws.Range("F25").Value = 0.9498469295364202
$ws.Range("G25").Value = 0.4608617679631664
$ws.Range("H25").Value = 0.5750480193390644
$ws.Range("I25").Value = 0.7086108388372634
$ws.Range("J25").Value = 0.2601537845990904
$ws.Range("L25").Value = 0.4240433117470559
$ws.Range("N25").Value = 1.417945913309737
$ws.Range("O25").Value = 2.027479561103377
